# fall 23 week 2 inputs
# Appends 20 new matchup rows (1123-1142) to sheet1, mirroring the prior
# week's data layout: Player_1, Points_1, Player_2, Points_2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(7,12,6,8),
    @(2,7,4,13),
    @(6,12,5,8),
    @(2,6,3,14),
    @(4,14,5,6),
    @(3,5,4,15),
    @(6,8,4,12),
    @(5,8,4,12),
    @(3,19,4,1),
    @(6,8,5,12),
    @(3,13,2,7),
    @(2,4,3,16),
    @(4,12,5,8),
    @(9,15,4,5),
    @(3,13,5,7),
    @(5,13,4,7),
    @(4,5,3,15),
    @(6,5,7,15),
    @(4,4,3,16),
    @(5,6,3,14)
)

$startRow = 1123
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$endRow = $startRow + $data.Length - 1

# Match the updated view/selection state from the diff.
$excel.ActiveWindow.ScrollRow = 1121
$ws.Range("A1143").Select()
